# Regenerate save_data column G ("K" = strike count) from the refreshed
# source computation. The sheet previously held the raw "Strike#" count;
# this writes the recalculated K values (std/mean derived s_vals) in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 2
    4  = 2
    5  = 4
    6  = 1
    7  = 2
    8  = 0
    9  = 1
    10 = 3
    11 = 0
    12 = 2
    13 = 3
    14 = 2
    15 = 0
    16 = 1
    17 = 0
    18 = 1
    19 = 0
    20 = 3
    21 = 1
    22 = 1
    23 = 3
    24 = 1
    25 = 0
    26 = 0
    27 = 1
    28 = 0
    29 = 2
    30 = 1
    31 = 1
    32 = 1
    33 = 2
    34 = 1
    35 = 1
    36 = 0
    37 = 2
    38 = 1
    39 = 0
    40 = 1
    41 = 1
    42 = 0
    43 = 1
    44 = 0
    45 = 2
    46 = 1
    47 = 2
    48 = 1
    49 = 0
    50 = 3
    51 = 1
    52 = 1
    53 = 2
    54 = 2
    55 = 1
    56 = 1
    57 = 1
    58 = 2
    59 = 1
    60 = 1
    64 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
